# Update the "Förändrad" (Changed) date column (C) for every data row,
# bumping each date forward by one day (45180 -> 45181, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 2 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $cur = $cell.Value2
    if ($cur -ne $null) {
        $cell.Value2 = $cur + 1
    }
}
